$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for semana 26 (row 27) and semana 28 (row 29)
$ws.Range("B27").Value = 272
$ws.Range("B29").Value = 304

# Add new row for semana 29 (week 29 of 2025)
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 1
